# -----------------------------------------------------------------------
# Helper: write a text value to a cell without Excel's COM layer
# reinterpreting it as a date/time/number (e.g. "2026-02-17") and without
# leaving any stray formatting (quote-prefix style) behind. Writing through
# Formula with a leading apostrophe forces literal text entry (exactly like
# a user typing '2026-02-17 into the formula bar), then ClearFormats()
# drops the quote-prefix style so the cell matches a plain default-styled
# text cell.
# -----------------------------------------------------------------------
function Set-SafeText {
    param($range, [string]$text)
    $range.Formula = "'" + $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Summary sheet: Total Trades (B6) 67 -> 68, Win Rate % (B9) 47.76 -> 47.06
# -----------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 68
$summary.Range("B9").Value = 47.06

# -----------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 5)
# Trades (D5) 34 -> 35, Win Rate % (G5) 52.94 -> 51.43
# -----------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 35
$status.Range("G5").Value = 51.43

# -----------------------------------------------------------------
# All Trades sheet: trade row 69 (Trade #68) goes from OPEN to CLOSED
# -----------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G69").Value = 0.057034
Set-SafeText $allTrades.Range("H69") "CLOSED"
$allTrades.Range("I69").Value = -4.943
$allTrades.Range("J69").Value = 0
$allTrades.Range("K69").Value = 100.58
Set-SafeText $allTrades.Range("L69") "early_exit"
$allTrades.Range("M69").Value = 0.14

# New trade row 102 (Trade #101) appended to All Trades
$allTrades.Range("A102").Value = 101
Set-SafeText $allTrades.Range("B102") "2026-02-17"
Set-SafeText $allTrades.Range("C102") "21:03:31"
Set-SafeText $allTrades.Range("D102") "MarketMaking"
Set-SafeText $allTrades.Range("E102") "UP"
$allTrades.Range("F102").Value = 0.06
Set-SafeText $allTrades.Range("G102") ""
Set-SafeText $allTrades.Range("H102") "OPEN"
$allTrades.Range("I102").Value = 0
$allTrades.Range("J102").Value = 0
$allTrades.Range("K102").Value = 100.584887765177
Set-SafeText $allTrades.Range("L102") ""
$allTrades.Range("M102").Value = 0
$allTrades.Range("N102").Value = 0
$allTrades.Range("O102").Value = 0
$allTrades.Range("P102").Value = 0.6
Set-SafeText $allTrades.Range("Q102") "Normal spread capture: 19600 bps"

# -----------------------------------------------------------------
# MarketMaking sheet: trade row 36 (Trade #68) goes from OPEN to CLOSED
# (column layout differs from All Trades sheet)
# -----------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G36").Value = 0.057034
Set-SafeText $mm.Range("H36") "CLOSED"
$mm.Range("I36").Value = -4.943
$mm.Range("J36").Value = 0
$mm.Range("K36").Value = 100.58
Set-SafeText $mm.Range("P36") "early_exit"
$mm.Range("Q36").Value = 0.14

# New trade row 69 (Trade #101) appended to MarketMaking
$mm.Range("A69").Value = 101
Set-SafeText $mm.Range("B69") "2026-02-17"
Set-SafeText $mm.Range("C69") "21:03:31"
Set-SafeText $mm.Range("D69") "MarketMaking"
Set-SafeText $mm.Range("E69") "UP"
$mm.Range("F69").Value = 0.06
Set-SafeText $mm.Range("G69") ""
Set-SafeText $mm.Range("H69") "OPEN"
$mm.Range("I69").Value = 0
$mm.Range("J69").Value = 0
$mm.Range("K69").Value = 100.584887765177
$mm.Range("L69").Value = 0
$mm.Range("M69").Value = 0
$mm.Range("N69").Value = 0.6
Set-SafeText $mm.Range("O69") "Normal spread capture: 19600 bps"
Set-SafeText $mm.Range("P69") ""
$mm.Range("Q69").Value = 0
